# Apply cryptos list update (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Worksheet, [string]$Address, [string]$Text)
    $range = $Worksheet.Range($Address)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = "Normal"
}

# Row 2
Set-TextCell $ws "D2" "67.769.97"
Set-TextCell $ws "E2" "  +0.00%  "

# Row 3
Set-TextCell $ws "D3" "3.779.84"
Set-TextCell $ws "E3" "  -0.66%  "

# Row 4
Set-TextCell $ws "D4" "0.999"
Set-TextCell $ws "E4" "  -0.10%  "

# Row 5
Set-TextCell $ws "D5" "599.07"
Set-TextCell $ws "E5" "  -0.09%  "

# Row 6
Set-TextCell $ws "D6" "162.47"
Set-TextCell $ws "E6" "  -3.08%  "

# Row 7
Set-TextCell $ws "D7" "3.776.82"
Set-TextCell $ws "E7" "  -0.75%  "

# Row 8
Set-TextCell $ws "E8" "  +0.09%  "

# Row 9
Set-TextCell $ws "E9" "  -1.45%  "

# Row 10
Set-TextCell $ws "E10" "  -3.16%  "

# Row 11
Set-TextCell $ws "D11" "0.445"
Set-TextCell $ws "E11" "  -1.07%  "

# Row 12
Set-TextCell $ws "D12" "6.60"
Set-TextCell $ws "E12" "  +4.78%  "

# Row 13
Set-TextCell $ws "E13" "  -3.88%  "

# Row 14
Set-TextCell $ws "D14" "35.01"
Set-TextCell $ws "E14" "  -2.69%  "

# Row 15
Set-TextCell $ws "D15" "4.411.04"
Set-TextCell $ws "E15" "  -0.68%  "

# Row 16
Set-TextCell $ws "D16" "3.771.63"
Set-TextCell $ws "E16" "  -1.55%  "

# Row 17
Set-TextCell $ws "D17" "67.752.48"

# Row 18
Set-TextCell $ws "D18" "18.11"
Set-TextCell $ws "E18" "  -2.01%  "

# Row 19
Set-TextCell $ws "E19" "  +1.67%  "

# Row 20
Set-TextCell $ws "E20" "  -1.26%  "

# Row 21
Set-TextCell $ws "D21" "455.88"
Set-TextCell $ws "E21" "  -1.34%  "

# Row 22
Set-TextCell $ws "D22" "9.42"
Set-TextCell $ws "E22" "  -4.84%  "

# Row 23
Set-TextCell $ws "D23" "0.690"
Set-TextCell $ws "E23" "  -1.50%  "

# Row 24
Set-TextCell $ws "E24" "  -0.69%  "

# Row 25
Set-TextCell $ws "E25" "  -6.69%  "

# Row 26
Set-TextCell $ws "D26" "11.83"
Set-TextCell $ws "E26" "  -2.16%  "

# Row 27
Set-TextCell $ws "B27" "Fetch.AI"
Set-TextCell $ws "C27" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell $ws "D27" "2.08"
Set-TextCell $ws "E27" "  -1.52%  "

# Row 28
Set-TextCell $ws "B28" "Dai"
Set-TextCell $ws "C28" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws "D28" "1.00"
Set-TextCell $ws "E28" "  -0.01%  "

# Row 29
Set-TextCell $ws "D29" "9.88"
Set-TextCell $ws "E29" "  -1.34%  "

# Row 30
Set-TextCell $ws "D30" "3.925.76"
Set-TextCell $ws "E30" "  -0.70%  "

# Row 31
Set-TextCell $ws "B31" "ImmutableX"
Set-TextCell $ws "C31" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws "D31" "2.20"
Set-TextCell $ws "E31" "  -2.03%  "

# Row 32
Set-TextCell $ws "B32" "NEARProtocol"
Set-TextCell $ws "C32" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws "D32" "7.19"
Set-TextCell $ws "E32" "  -3.57%  "

# Row 33
Set-TextCell $ws "E33" "  -7.00%  "

# Row 34
Set-TextCell $ws "D34" "28.83"
Set-TextCell $ws "E34" "  -2.37%  "

# Row 35
Set-TextCell $ws "E35" "  +0.31%  "

# Row 36
Set-TextCell $ws "D36" "8.91"

# Row 37
Set-TextCell $ws "D37" "0.0990"
Set-TextCell $ws "E37" "  -1.11%  "

# Row 38
Set-TextCell $ws "E38" "  +4.91%  "

# Row 39
Set-TextCell $ws "E39" "  -0.04%  "

# Row 40
Set-TextCell $ws "D40" "0.976"
Set-TextCell $ws "E40" "  -2.18%  "

# Row 41
Set-TextCell $ws "E41" "  -8.35%  "

# Row 42
Set-TextCell $ws "D42" "0.999"
Set-TextCell $ws "E42" "  -0.05%  "

# Row 43
Set-TextCell $ws "E43" "  +0.07%  "

# Row 44
Set-TextCell $ws "D44" "43.55"
Set-TextCell $ws "E44" "  +1.39%  "

# Row 45
Set-TextCell $ws "D45" "47.13"

# Row 46
Set-TextCell $ws "D46" "151.43"
Set-TextCell $ws "E46" "  +2.50%  "

# Row 47
Set-TextCell $ws "D47" "0.294"
Set-TextCell $ws "E47" "  -2.63%  "

# Row 48
Set-TextCell $ws "E48" "  -1.01%  "

# Row 49
Set-TextCell $ws "D49" "1.36"
Set-TextCell $ws "E49" "  -1.18%  "

# Row 51
Set-TextCell $ws "D51" "383.43"
Set-TextCell $ws "E51" "  -2.01%  "
